$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,23
$row2[0,0] = 0.979323308270677
$row2[0,1] = 0.0075187969924812
$row2[0,2] = 0.0093984962406015
$row2[0,3] = 0.0507518796992481
$row2[0,4] = 0.0056390977443609
$row2[0,5] = 0.0093984962406015
$row2[0,6] = 0.973684210526316
$row2[0,7] = 0.0244360902255639
$row2[0,8] = 0.951127819548872
$row2[0,9] = 0.964285714285714
$row2[0,10] = 0
$row2[0,11] = 0.0037593984962406
$row2[0,12] = 0.99812030075188
$row2[0,13] = 0
$row2[0,14] = 0.996240601503759
$row2[0,15] = 0
$row2[0,16] = 0.890977443609023
$row2[0,17] = 0.0526315789473684
$row2[0,18] = 0.0639097744360902
$row2[0,19] = 0.994360902255639
$row2[0,20] = 0.0056390977443609
$row2[0,21] = 0.0018796992481203
$row2[0,22] = 0.018796992481203
$ws.Range("B2:X2").Value = $row2

$row3 = New-Object "object[,]" 1,23
$row3[0,0] = 0.0056390977443609
$row3[0,1] = 0.0244360902255639
$row3[0,2] = 0.0037593984962406
$row3[0,3] = 0.0056390977443609
$row3[0,4] = 0.0037593984962406
$row3[0,5] = 0.977443609022556
$row3[0,6] = 0.0093984962406015
$row3[0,7] = 0.0131578947368421
$row3[0,8] = 0
$row3[0,9] = 0.0018796992481203
$row3[0,10] = 0.93609022556391
$row3[0,11] = 0.0056390977443609
$row3[0,12] = 0
$row3[0,13] = 0.99812030075188
$row3[0,14] = 0.0037593984962406
$row3[0,15] = 1
$row3[0,16] = 0.093984962406015
$row3[0,17] = 0.934210526315789
$row3[0,18] = 0.926691729323308
$row3[0,19] = 0
$row3[0,20] = 0.0037593984962406
$row3[0,21] = 0
$row3[0,22] = 0.0075187969924812
$ws.Range("B3:X3").Value = $row3

$row4 = New-Object "object[,]" 1,23
$row4[0,0] = 0.0037593984962406
$row4[0,1] = 0.0075187969924812
$row4[0,2] = 0
$row4[0,3] = 0.941729323308271
$row4[0,4] = 0.988721804511278
$row4[0,5] = 0.0112781954887218
$row4[0,6] = 0.0056390977443609
$row4[0,7] = 0.949248120300752
$row4[0,8] = 0.0093984962406015
$row4[0,9] = 0.0319548872180451
$row4[0,10] = 0.0018796992481203
$row4[0,11] = 0
$row4[0,12] = 0.0018796992481203
$row4[0,13] = 0
$row4[0,14] = 0
$row4[0,15] = 0
$row4[0,16] = 0.0112781954887218
$row4[0,17] = 0.0037593984962406
$row4[0,18] = 0.0018796992481203
$row4[0,19] = 0.0056390977443609
$row4[0,20] = 0.984962406015038
$row4[0,21] = 0.996240601503759
$row4[0,22] = 0.969924812030075
$ws.Range("B4:X4").Value = $row4

$row5 = New-Object "object[,]" 1,23
$row5[0,0] = 0.0112781954887218
$row5[0,1] = 0.960526315789474
$row5[0,2] = 0.986842105263158
$row5[0,3] = 0.0018796992481203
$row5[0,4] = 0.0018796992481203
$row5[0,5] = 0.0018796992481203
$row5[0,6] = 0.0112781954887218
$row5[0,7] = 0.0131578947368421
$row5[0,8] = 0.037593984962406
$row5[0,9] = 0
$row5[0,10] = 0.0601503759398496
$row5[0,11] = 0.990601503759398
$row5[0,12] = 0
$row5[0,13] = 0.0018796992481203
$row5[0,14] = 0
$row5[0,15] = 0
$row5[0,16] = 0.0037593984962406
$row5[0,17] = 0.0075187969924812
$row5[0,18] = 0.0075187969924812
$row5[0,19] = 0
$row5[0,20] = 0.0056390977443609
$row5[0,21] = 0.0018796992481203
$row5[0,22] = 0.0018796992481203
$ws.Range("B5:X5").Value = $row5

